$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the Kyrgyz title text in A1 (spelling correction).
$ws.Range("A1").Value = "3.9.2 Коопсуздук суунун, коопсуздук санитариянын жана гигиенанын жоктугунан болгон өлүм"

# 2. Add the new 2022 data column (S), mirroring the formatting of column R.
$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 1.2
$ws.Range("S6").Value = 2.7
$ws.Range("S7").Value = 0.9
$ws.Range("S8").Value = 0.4
$ws.Range("S9").Value = 0.7
$ws.Range("S10").Value = 0.9
$ws.Range("S11").Value = 1.1
$ws.Range("S12").Value = 2.7
$ws.Range("S13").Value = 0.4
$ws.Range("S14").Value = 0.6

$ws.Range("S4:S14").Style = "Обычный 6"

for ($r = 4; $r -le 14; $r++) {
    $src = $ws.Cells.Item($r, 18)
    $dst = $ws.Cells.Item($r, 19)
    $dst.Font.Name = $src.Font.Name
    $dst.Font.Size = $src.Font.Size
    $dst.Font.Bold = $src.Font.Bold
    $dst.NumberFormat = $src.NumberFormat
    $dst.HorizontalAlignment = $src.HorizontalAlignment
    $dst.VerticalAlignment = $src.VerticalAlignment
    $dst.Borders.Item(7).LineStyle = $src.Borders.Item(7).LineStyle
    $dst.Borders.Item(10).LineStyle = $src.Borders.Item(10).LineStyle
    $dst.Borders.Item(8).LineStyle = $src.Borders.Item(8).LineStyle
    $dst.Borders.Item(8).Weight = $src.Borders.Item(8).Weight
    $dst.Borders.Item(9).LineStyle = $src.Borders.Item(9).LineStyle
    $dst.Borders.Item(9).Weight = $src.Borders.Item(9).Weight
}
